# Updated cryptos list on Tue Oct 31 14:00:33 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto table on Sheet1, and re-syncs two pairs of rows whose coins swapped
# rank order (Polkadot/WrappedBTC at rows 16-17, and
# PancakeSwap/InternetComputer(DFINITY) at rows 30-31, and
# HuobiToken/MXToken/ARBITRUM at rows 41-43) along with their Coin name and
# Link columns.
#
# Several of the new Price values (e.g. "227.08", "0.640", "11.40") are
# strings that *look* numeric. Excel's COM layer auto-converts such strings
# to real numbers (and normalizes away meaningful trailing zeros) unless the
# cell is pre-formatted as Text ("@") before the value is assigned, so those
# cells get their NumberFormat set first to keep them as plain text, matching
# the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = '34.325.27'
$ws.Range('E2').Value = '  -1.12%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '1.799.42'
$ws.Range('E3').Value = '  -1.15%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  +0.34%  '

# Row 5 - BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.08'
$ws.Range('E5').Value = '  -0.68%  '

# Row 6 - XRP
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.598'
$ws.Range('E6').Value = '  +3.73%  '

# Row 7 - USDC
$ws.Range('E7').Value = '  +0.30%  '

# Row 8 - Solana
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '36.14'
$ws.Range('E8').Value = '  +3.91%  '

# Row 9 - Cardano
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.295'
$ws.Range('E9').Value = '  -1.91%  '

# Row 10 - Dogecoin
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0687'
$ws.Range('E10').Value = '  -1.68%  '

# Row 11 - TRON
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0964'
$ws.Range('E11').Value = '  +0.98%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '2.063.47'
$ws.Range('E12').Value = '  -0.92%  '

# Row 13 - Chainlink
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.24'
$ws.Range('E13').Value = '  -1.59%  '

# Row 14 - WrappedEther
$ws.Range('D14').Value = '1.799.39'
$ws.Range('E14').Value = '  -1.25%  '

# Row 15 - Polygon
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.640'
$ws.Range('E15').Value = '  -0.79%  '

# Row 16 - was Polkadot, now WrappedBTC
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '34.332.21'
$ws.Range('E16').Value = '  -1.13%  '

# Row 17 - was WrappedBTC, now Polkadot
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.42'
$ws.Range('E17').Value = '  +1.88%  '

# Row 18 - Litecoin
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.71'
$ws.Range('E18').Value = '  +0.90%  '

# Row 19 - BitcoinCash
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.68'
$ws.Range('E19').Value = '  -1.55%  '

# Row 20 - ShibaInu
$ws.Range('D20').Value = '0.0₃0782'
$ws.Range('E20').Value = '  -2.56%  '

# Row 21 - Avalanche
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.40'
$ws.Range('E21').Value = '  -1.78%  '

# Row 22 - Dai
$ws.Range('E22').Value = '  +0.30%  '

# Row 23 - Uniswap
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.14'
$ws.Range('E23').Value = '  -1.38%  '

# Row 24 - Toncoin
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.23'
$ws.Range('E24').Value = '  +6.18%  '

# Row 25 - Monero
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '169.87'
$ws.Range('E25').Value = '  -1.89%  '

# Row 26 - Cosmos
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.15'
$ws.Range('E26').Value = '  +8.97%  '

# Row 27 - EthereumClassic
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.27'
$ws.Range('E27').Value = '  +2.66%  '

# Row 28 - Stellar
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.121'
$ws.Range('E28').Value = '  +1.87%  '

# Row 29 - BinanceUSD
$ws.Range('E29').Value = '  +0.30%  '

# Row 30 - was InternetComputer(DFINITY), now PancakeSwap
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.24'
$ws.Range('E30').Value = '  -0.86%  '

# Row 31 - was PancakeSwap, now InternetComputer(DFINITY)
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.94'
$ws.Range('E31').Value = '  -1.59%  '

# Row 32 - Filecoin
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.80'
$ws.Range('E32').Value = '  -1.51%  '

# Row 33 - Hedera
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0518'
$ws.Range('E33').Value = '  -2.81%  '

# Row 34 - LidoDAOToken
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.78'
$ws.Range('E34').Value = '  -4.21%  '

# Row 35 - Maker
$ws.Range('D35').Value = '1.372.17'
$ws.Range('E35').Value = '  -3.36%  '

# Row 36 - ImmutableX
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.653'
$ws.Range('E36').Value = '  -4.85%  '

# Row 37 - TrustWalletToken
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.05'
$ws.Range('E37').Value = '  -1.87%  '

# Row 38 - RenderToken
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.34'
$ws.Range('E38').Value = '  -11.58%  '

# Row 39 - VeChain
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0186'
$ws.Range('E39').Value = '  -3.03%  '

# Row 40 - Aave
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '81.74'
$ws.Range('E40').Value = '  -4.16%  '

# Row 41 - was ARBITRUM, now HuobiToken
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.42'
$ws.Range('E41').Value = '  +0.47%  '

# Row 42 - was HuobiToken, now MXToken
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.79'
$ws.Range('E42').Value = '  -2.44%  '

# Row 43 - was MXToken, now ARBITRUM
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.944'
$ws.Range('E43').Value = '  -1.37%  '

# Row 45 - InjectiveProtocol
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.23'
$ws.Range('E45').Value = '  -4.54%  '

# Row 46 - Kaspa
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0497'
$ws.Range('E46').Value = '  -4.09%  '

# Row 47 - FraxShare
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.92'
$ws.Range('E47').Value = '  -3.22%  '

# Row 48 - RocketPoolETH
$ws.Range('D48').Value = '1.962.81'
$ws.Range('E48').Value = '  -1.05%  '

# Row 49 - PaxDollar
$ws.Range('E49').Value = '  +0.33%  '

# Row 50 - Quant
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.39'
$ws.Range('E50').Value = '  -3.04%  '

# Row 51 - BabyDogeCoin
$ws.Range('D51').Value = '0.0₆0121'
$ws.Range('E51').Value = '  -7.56%  '
